$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.615.43"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.53"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.64"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.02"
$ws.Range("E6").Value = "  +7.75%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +3.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.48"
$ws.Range("E10").Value = "  +12.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.649.71"
$ws.Range("E14").Value = "  +1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +3.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.328.63"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.808"
$ws.Range("E17").Value = "  +5.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.516.84"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +1.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  +1.91%  "

$ws.Range("E21").Value = "  +2.26%  "

$ws.Range("E22").Value = "  +1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.57"
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.07"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.14"
$ws.Range("E28").Value = "  +10.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.80"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.32"
$ws.Range("E32").Value = "  +0.97%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +5.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0753"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.39"
$ws.Range("E36").Value = "  +2.58%  "

$ws.Range("E37").Value = "  +3.78%  "

$ws.Range("E38").Value = "  +4.72%  "

$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("E41").Value = "  +7.48%  "

$ws.Range("E42").Value = "  +16.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.007.88"
$ws.Range("E43").Value = "  -1.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.28"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +3.68%  "

$ws.Range("E46").Value = "  +6.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.23"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.91"
$ws.Range("E48").Value = "  +4.18%  "

$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.02"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("E51").Value = "  -0.09%  "
